# Performance Tracking.xlsx — add low pass filter standard-deviation table
# (columns U:Z) and the accompanying "Filter disabled" baseline label that
# replaces the old duplicate Lane1-4/Mean shared strings, per the commit
# "Added low pass filter data to spreadsheet".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New table header (row 1): title spanning the new block
# ---------------------------------------------------------------------
$ws.Range("U1").Value = "Standard deviations for low pass filter parameters (using network 1503559070.h5)"

# ---------------------------------------------------------------------
# New table column headers (row 2)
# ---------------------------------------------------------------------
$ws.Range("U2").Value = "Parameters"
$ws.Range("V2").Value = "Lane 1"
$ws.Range("W2").Value = "Lane 2"
$ws.Range("X2").Value = "Lane 3"
$ws.Range("Y2").Value = "Lane 4"
$ws.Range("Z2").Value = "Mean"

# ---------------------------------------------------------------------
# Row 3: baseline "Filter disabled" values (same measurements that were
# already recorded for the 0915d0e / 1503559070.h5 run in O8:S8)
# ---------------------------------------------------------------------
$ws.Range("U3").Value = "Filter disabled"
$ws.Range("V3").Value = 0.7085852
$ws.Range("W3").Value = 0.3898657
$ws.Range("X3").Value = 0.6391341
$ws.Range("Y3").Value = 0.8785638
$ws.Range("Z3").Formula = "=AVERAGE(V3:Y3)"

# ---------------------------------------------------------------------
# Row 4: 1.0, 0.3, 0.2, 0.1
# ---------------------------------------------------------------------
$ws.Range("U4").Value = "1.0, 0.3, 0.2, 0.1"
$ws.Range("V4").Value = 0.7065313
$ws.Range("W4").Value = 0.3907571
$ws.Range("X4").Value = 0.6443396
$ws.Range("Y4").Value = 0.9180865
$ws.Range("Z4").Formula = "=AVERAGE(V4:Y4)"

# ---------------------------------------------------------------------
# Row 5: 1.0, 0.7, 0.6, 0.5, 0.4
# ---------------------------------------------------------------------
$ws.Range("U5").Value = "1.0, 0.7, 0.6, 0.5, 0.4"
$ws.Range("V5").Value = 0.6531526
$ws.Range("W5").Value = 0.3540962
$ws.Range("X5").Value = 0.6581528
$ws.Range("Y5").Value = 0.8799137
$ws.Range("Z5").Formula = "=AVERAGE(V5:Y5)"

# ---------------------------------------------------------------------
# Row 6: 1.0, 0.8, 0.3, 0.1
# ---------------------------------------------------------------------
$ws.Range("U6").Value = "1.0, 0.8, 0.3, 0.1"
$ws.Range("V6").Value = 0.6268617
$ws.Range("W6").Value = 1.001549
$ws.Range("X6").Value = 0.650547
$ws.Range("Y6").Value = 0.8531444
$ws.Range("Z6").Formula = "=AVERAGE(V6:Y6)"

# ---------------------------------------------------------------------
# Touch the previously-empty cells across columns A:S so each gets the
# workbook's default ("Normal") style instead of inheriting the bold
# column style - mirrors the blank, explicitly-styled cells introduced
# around the existing tables.
# ---------------------------------------------------------------------
$blankCells = @(
    "B1","C1","D1","E1","F1","G1","H1","I1",
    "L1","M1","N1","O1","P1","Q1","R1","S1",
    "B2","C2","E2","F2","G2","H2","I2",
    "L2","N2","P2","Q2","R2","S2",
    "P4","Q4","R4","S4",
    "P6","Q6","R6","S6",
    "P7","Q7","R7","S7"
)
foreach ($cell in $blankCells) {
    $ws.Range($cell).Font.Size = 10
}

# ---------------------------------------------------------------------
# Restore the selection to the cell the author ended up on (B12).
# ---------------------------------------------------------------------
$ws.Range("B12").Select()
